$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Refresh the "datetime1" auto-date placeholders on the slide master
#    and every slide layout: 11/2/2022 -> 9/27/2023 (the deck was
#    re-saved on 9/27/2023, which is when PowerPoint recalculates the
#    cached text of every auto-updating date field in the masters).
# ---------------------------------------------------------------------
$newDate = "9/27/2023"

function Set-DatePlaceholderText {
    param($shapes)
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide master
Set-DatePlaceholderText $p.SlideMaster.Shapes

# Every slide layout used by the master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Set-DatePlaceholderText $layout.Shapes
}

# ---------------------------------------------------------------------
# 2) Slide 1 subtitle: drop the "Fall 2022 | " prefix and nudge the
#    textbox position to match the new, shorter text box.
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$subtitle = $s1.Shapes.Item(3)
$subtitle.TextFrame.TextRange.Text = "University of Mount union"
# Shape.Left/.Top are single-precision (points); nudge by the smallest
# representable step so the EMU the file round-trips to lands on the
# exact target (581191, 1440465) instead of its float32 neighbour.
$subtitle.Left = 581192 / 12700
$subtitle.Top = 1440465 / 12700
